# Update country data and fix row order (Nepal/Costa Rica, Afganistan/Azerbaiyan,
# Islas Malvinas/Montserrat) plus refresh the "datos actualizados" timestamp,
# as part of updating countries & provincias Spain COVID data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = 'Datos actualizados a 18 de Septiembre de 2020 a las 14:04'
$ws.Range("B4").Value = 6876126
$ws.Range("C4").Value = 1530
$ws.Range("D4").Value = 4155933
$ws.Range("E4").Value = 2517956
$ws.Range("G4").Value = 24
$ws.Range("H4").Value = 202237
$ws.Range("B31").Value = 122917
$ws.Range("C31").Value = 224
$ws.Range("D31").Value = 119822
$ws.Range("E31").Value = 2886
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 209
$ws.Range("A55").Value = 'Nepal'
$ws.Range("B55").Value = 61593
$ws.Range("C55").Value = 2020
$ws.Range("D55").Value = 43820
$ws.Range("E55").Value = 17383
$ws.Range("G55").Value = 7
$ws.Range("H55").Value = 390
$ws.Range("A56").Value = 'Costa Rica'
$ws.Range("B56").Value = 60818
$ws.Range("D56").Value = 22662
$ws.Range("E56").Value = 37490
$ws.Range("H56").Value = 666
$ws.Range("A67").Value = 'Azerbaiyan'
$ws.Range("B67").Value = 38894
$ws.Range("C67").Value = 117
$ws.Range("D67").Value = 36424
$ws.Range("E67").Value = 1898
$ws.Range("H67").Value = 572
$ws.Range("A68").Value = 'Afganistan'
$ws.Range("B68").Value = 38883
$ws.Range("C68").Value = 11
$ws.Range("D68").Value = 32576
$ws.Range("E68").Value = 4870
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 1437
$ws.Range("B79").Value = 24897
$ws.Range("C79").Value = 292
$ws.Range("D79").Value = 17365
$ws.Range("E79").Value = 6780
$ws.Range("G79").Value = 5
$ws.Range("H79").Value = 752
$ws.Range("B81").Value = 21847
$ws.Range("C81").Value = 454
$ws.Range("D81").Value = 17110
$ws.Range("E81").Value = 4102
$ws.Range("B88").Value = 14645
$ws.Range("C88").Value = 27
$ws.Range("D88").Value = 11051
$ws.Range("E88").Value = 3293
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 301
$ws.Range("B150").Value = 2230
$ws.Range("C150").Value = 24
$ws.Range("E150").Value = 108
$ws.Range("B179").Value = 430
$ws.Range("C179").Value = 1
$ws.Range("E179").Value = 18
$ws.Range("B182").Value = 346
$ws.Range("C182").Value = 3
$ws.Range("D182").Value = 315
$ws.Range("E182").Value = 31
$ws.Range("A214").Value = 'Montserrat'
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
$ws.Range("A215").Value = 'Islas Malvinas'
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
